$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$bf2 = New-Object 'object[,]' 1,5
$bf2[0,0] = 1.02
$bf2[0,1] = 1.032457870620217
$bf2[0,2] = 1.041018084730279
$bf2[0,3] = 1.050467848892552
$bf2[0,4] = 1.055279202385765
$ws.Range("B2:F2").Value = $bf2

$inb2 = New-Object 'object[,]' 1,6
$inb2[0,0] = 1.03580434566463
$inb2[0,1] = 1.03758798229482
$inb2[0,2] = 1.043798789197469
$inb2[0,3] = 1.053222043623242
$inb2[0,4] = 1.058020097333649
$inb2[0,5] = 1.039061476242684
$ws.Range("I2:N2").Value = $inb2

$bf3 = New-Object 'object[,]' 1,5
$bf3[0,0] = 1.02
$bf3[0,1] = 1.033368681826004
$bf3[0,2] = 1.041700024947126
$bf3[0,3] = 1.051344080027163
$bf3[0,4] = 1.056123803695223
$ws.Range("B3:F3").Value = $bf3

$inb3 = New-Object 'object[,]' 1,6
$inb3[0,0] = 1.035944968108676
$inb3[0,1] = 1.038141293066593
$inb3[0,2] = 1.044291609357001
$inb3[0,3] = 1.053910562723084
$inb3[0,4] = 1.058678028897196
$inb3[0,5] = 1.039615572779218
$ws.Range("I3:N3").Value = $inb3

$bf4 = New-Object 'object[,]' 1,5
$bf4[0,0] = 1.02
$bf4[0,1] = 1.033958616341583
$bf4[0,2] = 1.042141658999556
$bf4[0,3] = 1.051912287537551
$bf4[0,4] = 1.056671193991466
$ws.Range("B4:F4").Value = $bf4

$inb4 = New-Object 'object[,]' 1,6
$inb4[0,0] = 1.036034794726805
$inb4[0,1] = 1.038499274299852
$inb4[0,2] = 1.044610190651168
$inb4[0,3] = 1.054356666700077
$inb4[0,4] = 1.059103983300465
$inb4[0,5] = 1.039974062386889
$ws.Range("I4:N4").Value = $inb4

$bf5 = New-Object 'object[,]' 1,5
$bf5[0,0] = 1.02
$bf5[0,1] = 1.034206762094523
$bf5[0,2] = 1.042327409653594
$bf5[0,3] = 1.052151454063709
$bf5[0,4] = 1.056901525016691
$ws.Range("B5:F5").Value = $bf5

$inb5 = New-Object 'object[,]' 1,6
$inb5[0,0] = 1.036072278196268
$inb5[0,1] = 1.038649757016874
$inb5[0,2] = 1.044744047688809
$inb5[0,3] = 1.054544347661428
$inb5[0,4] = 1.059283107818413
$inb5[0,5] = 1.04012475880663
$ws.Range("I5:N5").Value = $inb5

$bf6 = New-Object 'object[,]' 1,5
$bf6[0,0] = 1.02
$bf6[0,1] = 1.034248434871711
$bf6[0,2] = 1.042358603101004
$bf6[0,3] = 1.052191628253281
$bf6[0,4] = 1.056940210763981
$ws.Range("B6:F6").Value = $bf6

$inb6 = New-Object 'object[,]' 1,6
$inb6[0,0] = 1.036078555419564
$inb6[0,1] = 1.038675022955561
$inb6[0,2] = 1.044766518468318
$inb6[0,3] = 1.054575868220461
$inb6[0,4] = 1.059313186695451
$inb6[0,5] = 1.040150060625848
$ws.Range("I6:N6").Value = $inb6

$bf7 = New-Object 'object[,]' 1,5
$bf7[0,0] = 1.02
$bf7[0,1] = 1.033961931539439
$bf7[0,2] = 1.042144140666342
$bf7[0,3] = 1.051915482146863
$bf7[0,4] = 1.056674270870237
$ws.Range("B7:F7").Value = $bf7

$inb7 = New-Object 'object[,]' 1,6
$inb7[0,0] = 1.036035296682324
$inb7[0,1] = 1.038501285108196
$inb7[0,2] = 1.044611979549859
$inb7[0,3] = 1.054359173958517
$inb7[0,4] = 1.059106376563494
$inb7[0,5] = 1.039976076050812
$ws.Range("I7:N7").Value = $inb7

$bf8 = New-Object 'object[,]' 1,5
$bf8[0,0] = 1.02
$bf8[0,1] = 1.032765562928059
$bf8[0,2] = 1.041248471380206
$bf8[0,3] = 1.050763720120359
$bf8[0,4] = 1.055564456755162
$ws.Range("B8:F8").Value = $bf8

$inb8 = New-Object 'object[,]' 1,6
$inb8[0,0] = 1.035852110632947
$inb8[0,1] = 1.037774985478631
$inb8[0,2] = 1.043965402772963
$inb8[0,3] = 1.053454609647551
$inb8[0,4] = 1.058242400156061
$inb8[0,5] = 1.039248744992465
$ws.Range("I8:N8").Value = $inb8

$bf9 = New-Object 'object[,]' 1,5
$bf9[0,0] = 1.02
$bf9[0,1] = 1.030661897570593
$bf9[0,2] = 1.039673117415286
$bf9[0,3] = 1.048743645768536
$bf9[0,4] = 1.053615615166834
$ws.Range("B9:F9").Value = $bf9

$inb9 = New-Object 'object[,]' 1,6
$inb9[0,0] = 1.035520414288817
$inb9[0,1] = 1.036494836831787
$inb9[0,2] = 1.042823759385441
$inb9[0,3] = 1.051865204954897
$inb9[0,4] = 1.056721779414803
$inb9[0,5] = 1.037966778387707
$ws.Range("I9:N9").Value = $inb9

$bf10 = New-Object 'object[,]' 1,5
$bf10[0,0] = 1.02
$bf10[0,1] = 1.029262547370349
$bf10[0,2] = 1.038624948155136
$bf10[0,3] = 1.047403405722586
$bf10[0,4] = 1.052321058958289
$ws.Range("B10:F10").Value = $bf10

$inb10 = New-Object 'object[,]' 1,6
$inb10[0,0] = 1.035293336788262
$inb10[0,1] = 1.035641259690165
$inb10[0,2] = 1.042061194240933
$inb10[0,3] = 1.050808749989494
$inb10[0,4] = 1.0557093401832
$inb10[0,5] = 1.037111989068638
$ws.Range("I10:N10").Value = $inb10

$bf11 = New-Object 'object[,]' 1,5
$bf11[0,0] = 1.02
$bf11[0,1] = 1.028657362598425
$bf11[0,2] = 1.038171589881636
$bf11[0,3] = 1.046824624965578
$bf11[0,4] = 1.051761631834001
$ws.Range("B11:F11").Value = $bf11

$inb11 = New-Object 'object[,]' 1,6
$inb11[0,0] = 1.035193606950234
$inb11[0,1] = 1.035271631200706
$inb11[0,2] = 1.041730661616871
$inb11[0,3] = 1.050352058300362
$inb11[0,4] = 1.055271271303099
$inb11[0,5] = 1.036741835664327
$ws.Range("I11:N11").Value = $inb11

$bf12 = New-Object 'object[,]' 1,5
$bf12[0,0] = 1.02
$bf12[0,1] = 1.028432682895183
$bf12[0,2] = 1.038003270082852
$bf12[0,3] = 1.046609874967013
$bf12[0,4] = 1.051554006368268
$ws.Range("B12:F12").Value = $bf12

$inb12 = New-Object 'object[,]' 1,6
$inb12[0,0] = 1.035156352543513
$inb12[0,1] = 1.035134332246225
$inb12[0,2] = 1.041607837668028
$inb12[0,3] = 1.050182538365838
$inb12[0,4] = 1.055108603153603
$inb12[0,5] = 1.036604341729582
$ws.Range("I12:N12").Value = $inb12

$bf13 = New-Object 'object[,]' 1,5
$bf13[0,0] = 1.02
$bf13[0,1] = 1.028480872327043
$bf13[0,2] = 1.03803937171095
$bf13[0,3] = 1.046655928910516
$bf13[0,4] = 1.05159853498357
$ws.Range("B13:F13").Value = $bf13

$inb13 = New-Object 'object[,]' 1,6
$inb13[0,0] = 1.035164353247321
$inb13[0,1] = 1.035163783432026
$inb13[0,2] = 1.041634186038903
$inb13[0,3] = 1.050218895708295
$inb13[0,4] = 1.055143493731548
$inb13[0,5] = 1.036633834739444
$ws.Range("I13:N13").Value = $inb13

$bf14 = New-Object 'object[,]' 1,5
$bf14[0,0] = 1.02
$bf14[0,1] = 1.028638788194515
$bf14[0,2] = 1.038157674921056
$bf14[0,3] = 1.046806868860663
$bf14[0,4] = 1.051744465956246
$ws.Range("B14:F14").Value = $bf14

$inb14 = New-Object 'object[,]' 1,6
$inb14[0,0] = 1.035190531776036
$inb14[0,1] = 1.035260282070304
$inb14[0,2] = 1.041720509952074
$inb14[0,3] = 1.050338043361092
$inb14[0,4] = 1.055257824071869
$inb14[0,5] = 1.036730470416858
$ws.Range("I14:N14").Value = $inb14

$bf15 = New-Object 'object[,]' 1,5
$bf15[0,0] = 1.02
$bf15[0,1] = 1.028736100365026
$bf15[0,2] = 1.038230575771649
$bf15[0,3] = 1.046899899131023
$bf15[0,4] = 1.051834401514105
$ws.Range("B15:F15").Value = $bf15

$inb15 = New-Object 'object[,]' 1,6
$inb15[0,0] = 1.035206633382115
$inb15[0,1] = 1.035319737771752
$inb15[0,2] = 1.041773690451853
$inb15[0,3] = 1.050411469536752
$inb15[0,4] = 1.055328273455375
$inb15[0,5] = 1.036790010552222
$ws.Range("I15:N15").Value = $inb15

$bf16 = New-Object 'object[,]' 1,5
$bf16[0,0] = 1.02
$bf16[0,1] = 1.0293027272512
$bf16[0,2] = 1.038655046850021
$bf16[0,3] = 1.047441850345786
$bf16[0,4] = 1.052358210131966
$ws.Range("B16:F16").Value = $bf16

$inb16 = New-Object 'object[,]' 1,6
$inb16[0,0] = 1.035299926004098
$inb16[0,1] = 1.035665790306819
$inb16[0,2] = 1.042083123607078
$inb16[0,3] = 1.050839075256496
$inb16[0,4] = 1.055738420322872
$inb16[0,5] = 1.037136554521582
$ws.Range("I16:N16").Value = $inb16

$bf17 = New-Object 'object[,]' 1,5
$bf17[0,0] = 1.02
$bf17[0,1] = 1.029658357138572
$bf17[0,2] = 1.038921442953308
$bf17[0,3] = 1.047782218795926
$bf17[0,4] = 1.052687083683927
$ws.Range("B17:F17").Value = $bf17

$inb17 = New-Object 'object[,]' 1,6
$inb17[0,0] = 1.035358070590523
$inb17[0,1] = 1.035882854298735
$inb17[0,2] = 1.042277133502203
$inb17[0,3] = 1.05110750561347
$inb17[0,4] = 1.055995782389703
$inb17[0,5] = 1.037353926769263
$ws.Range("I17:N17").Value = $inb17

$bf18 = New-Object 'object[,]' 1,5
$bf18[0,0] = 1.02
$bf18[0,1] = 1.029865861650989
$bf18[0,2] = 1.039076875998653
$bf18[0,3] = 1.047980899583678
$bf18[0,4] = 1.052879018418458
$ws.Range("B18:F18").Value = $bf18

$inb18 = New-Object 'object[,]' 1,6
$inb18[0,0] = 1.035391849859291
$inb18[0,1] = 1.036009461666161
$inb18[0,2] = 1.042390263606844
$inb18[0,3] = 1.051264149756879
$inb18[0,4] = 1.056145928491151
$inb18[0,5] = 1.037480713933673
$ws.Range("I18:N18").Value = $inb18

$bf19 = New-Object 'object[,]' 1,5
$bf19[0,0] = 1.02
$bf19[0,1] = 1.029936627437437
$bf19[0,2] = 1.039129882861752
$bf19[0,3] = 1.048048669923816
$bf19[0,4] = 1.052944481528074
$ws.Range("B19:F19").Value = $bf19

$inb19 = New-Object 'object[,]' 1,6
$inb19[0,0] = 1.035403344718269
$inb19[0,1] = 1.03605263106766
$inb19[0,2] = 1.042428832485601
$inb19[0,3] = 1.051317573742818
$inb19[0,4] = 1.056197129672782
$inb19[0,5] = 1.037523944640674
$ws.Range("I19:N19").Value = $inb19

$bf20 = New-Object 'object[,]' 1,5
$bf20[0,0] = 1.02
$bf20[0,1] = 1.029620194000876
$bf20[0,2] = 1.038892856122911
$bf20[0,3] = 1.047745684992795
$bf20[0,4] = 1.052651787443704
$ws.Range("B20:F20").Value = $bf20

$inb20 = New-Object 'object[,]' 1,6
$inb20[0,0] = 1.035351846236183
$inb20[0,1] = 1.035859565625418
$inb20[0,2] = 1.042256321438292
$inb20[0,3] = 1.051078697988574
$inb20[0,4] = 1.055968166641782
$inb20[0,5] = 1.037330605023359
$ws.Range("I20:N20").Value = $inb20

$bf21 = New-Object 'object[,]' 1,5
$bf21[0,0] = 1.02
$bf21[0,1] = 1.028592282781956
$bf21[0,2] = 1.038122835417641
$bf21[0,3] = 1.046762414307579
$bf21[0,4] = 1.051701488196299
$ws.Range("B21:F21").Value = $bf21

$inb21 = New-Object 'object[,]' 1,6
$inb21[0,0] = 1.035182828652434
$inb21[0,1] = 1.035231865685815
$inb21[0,2] = 1.041695091065696
$inb21[0,3] = 1.050302954141526
$inb21[0,4] = 1.055224155240466
$inb21[0,5] = 1.036702013677845
$ws.Range("I21:N21").Value = $inb21

$bf22 = New-Object 'object[,]' 1,5
$bf22[0,0] = 1.02
$bf22[0,1] = 1.027946647823908
$bf22[0,2] = 1.037639142906059
$bf22[0,3] = 1.046145553196707
$bf22[0,4] = 1.051104985701293
$ws.Range("B22:F22").Value = $bf22

$inb22 = New-Object 'object[,]' 1,6
$inb22[0,0] = 1.035075344129823
$inb22[0,1] = 1.034837192053038
$inb22[0,2] = 1.041341937946712
$inb22[0,3] = 1.049815883287485
$inb22[0,4] = 1.054756656547499
$inb22[0,5] = 1.036306779563238
$ws.Range("I22:N22").Value = $inb22

$bf23 = New-Object 'object[,]' 1,5
$bf23[0,0] = 1.02
$bf23[0,1] = 1.028288848764843
$bf23[0,2] = 1.037895514250088
$bf23[0,3] = 1.046472433420977
$bf23[0,4] = 1.051421108684046
$ws.Range("B23:F23").Value = $bf23

$inb23 = New-Object 'object[,]' 1,6
$inb23[0,0] = 1.035132438815198
$inb23[0,1] = 1.035046416942946
$inb23[0,2] = 1.041529177740895
$inb23[0,3] = 1.050074024755214
$inb23[0,4] = 1.055004458396743
$inb23[0,5] = 1.036516301576488
$ws.Range("I23:N23").Value = $inb23

$bf24 = New-Object 'object[,]' 1,5
$bf24[0,0] = 1.02
$bf24[0,1] = 1.029637438058918
$bf24[0,2] = 1.038905773130765
$bf24[0,3] = 1.047762192583734
$bf24[0,4] = 1.052667735960759
$ws.Range("B24:F24").Value = $bf24

$inb24 = New-Object 'object[,]' 1,6
$inb24[0,0] = 1.035354659172676
$inb24[0,1] = 1.035870088781935
$inb24[0,2] = 1.042265725614943
$inb24[0,3] = 1.051091714688074
$inb24[0,4] = 1.055980644913443
$inb24[0,5] = 1.037341143123965
$ws.Range("I24:N24").Value = $inb24

$bf25 = New-Object 'object[,]' 1,5
$bf25[0,0] = 1.02
$bf25[0,1] = 1.031205205786848
$bf25[0,2] = 1.04008002699369
$bf25[0,3] = 1.049264749557733
$bf25[0,4] = 1.054118621696349
$ws.Range("B25:F25").Value = $bf25

$inb25 = New-Object 'object[,]' 1,6
$inb25[0,0] = 1.035607216563466
$inb25[0,1] = 1.036825816407192
$inb25[0,2] = 1.043119165842108
$inb25[0,3] = 1.052275555641649
$inb25[0,4] = 1.05711467215876
$inb25[0,5] = 1.038298227992074
$ws.Range("I25:N25").Value = $inb25
